$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Payment Methods")

# Insert a new row above the current row 2 (DBTRACCTID), shifting everything
# else down by one, then populate it with the "Check Form Code" entry that
# moved into the cash-account-details sheet.
$ws.Rows("2:2").Insert()

$ws.Range("A2").Value = "CHQFRMSCD"
$ws.Range("B2").Value = "Check Form Code"
$ws.Range("C2").Value = "'True"
$ws.Range("D2").Value = 18
$ws.Range("F2").Value = "^.{1,35}$"

# Copy the style used by the other data rows (row 3, formerly row 2) onto
# the newly inserted row so formatting matches.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# The Debtor ID row's "Required" flag flips to False now that the form code
# row has been added ahead of it. Force text (not a boolean) the same way
# Excel does when you type a value preceded by an apostrophe, then restore
# the plain (non quote-prefixed) cell format from an untouched neighbour.
$ws.Range("C6").Value = "'False"
$ws.Range("C7").Copy()
$ws.Range("C6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
